$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 3144
    5  = 19
    6  = 1756
    8  = 101
    10 = 9
    11 = 1451
    13 = 572
    14 = 360
    15 = 83
    17 = 82
    19 = 6
    23 = 3420
    24 = 412
    25 = 297
    26 = 470
    27 = 68
    28 = 20
    30 = 1153
}

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
